# Updates the 100 arithmetic-expression cells in the single 20x5 table
# (each <w:t> run holding an expression like "86-82=") to the values
# from the "commit" (master re-generated at 60844e3). Cells are addressed
# by (row, col) rather than Find/Replace on text, because a couple of the
# old expressions are duplicated across cells and Find scoped to a single
# cell.Range can still match/replace text in a different cell of the same
# table in this runtime; setting Range.Text directly is scoped correctly.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row 1, col 1: "86-82=" -> "18+76="
$t.Cell(1, 1).Range.Text = "18+76="
# row 1, col 2: "0+33=" -> "8+9="
$t.Cell(1, 2).Range.Text = "8+9="
# row 1, col 3: "11+83=" -> "46-5="
$t.Cell(1, 3).Range.Text = "46-5="
# row 1, col 4: "18-10=" -> "63+36="
$t.Cell(1, 4).Range.Text = "63+36="
# row 1, col 5: "5+59=" -> "47-11="
$t.Cell(1, 5).Range.Text = "47-11="
# row 2, col 1: "37+52=" -> "14+20="
$t.Cell(2, 1).Range.Text = "14+20="
# row 2, col 2: "16-5=" -> "92-18="
$t.Cell(2, 2).Range.Text = "92-18="
# row 2, col 3: "52-25=" -> "54-32="
$t.Cell(2, 3).Range.Text = "54-32="
# row 2, col 4: "14+49=" -> "36+57="
$t.Cell(2, 4).Range.Text = "36+57="
# row 2, col 5: "0+19=" -> "22+49="
$t.Cell(2, 5).Range.Text = "22+49="
# row 3, col 1: "64-8=" -> "96-47="
$t.Cell(3, 1).Range.Text = "96-47="
# row 3, col 2: "68-8=" -> "50+48="
$t.Cell(3, 2).Range.Text = "50+48="
# row 3, col 3: "1+37=" -> "30+53="
$t.Cell(3, 3).Range.Text = "30+53="
# row 3, col 4: "1+91=" -> "80+0="
$t.Cell(3, 4).Range.Text = "80+0="
# row 3, col 5: "48-42=" -> "6+53="
$t.Cell(3, 5).Range.Text = "6+53="
# row 4, col 1: "97-5=" -> "7+15="
$t.Cell(4, 1).Range.Text = "7+15="
# row 4, col 2: "18+78=" -> "74-56="
$t.Cell(4, 2).Range.Text = "74-56="
# row 4, col 3: "16+82=" -> "80-9="
$t.Cell(4, 3).Range.Text = "80-9="
# row 4, col 4: "29-18=" -> "8+34="
$t.Cell(4, 4).Range.Text = "8+34="
# row 4, col 5: "66-48=" -> "87-9="
$t.Cell(4, 5).Range.Text = "87-9="
# row 5, col 1: "92-20=" -> "60+18="
$t.Cell(5, 1).Range.Text = "60+18="
# row 5, col 2: "54-15=" -> "29+65="
$t.Cell(5, 2).Range.Text = "29+65="
# row 5, col 3: "88-77=" -> "82-60="
$t.Cell(5, 3).Range.Text = "82-60="
# row 5, col 4: "7+22=" -> "35-27="
$t.Cell(5, 4).Range.Text = "35-27="
# row 5, col 5: "84-42=" -> "80-47="
$t.Cell(5, 5).Range.Text = "80-47="
# row 6, col 1: "50-49=" -> "68+9="
$t.Cell(6, 1).Range.Text = "68+9="
# row 6, col 2: "87+8=" -> "87-36="
$t.Cell(6, 2).Range.Text = "87-36="
# row 6, col 3: "85+1=" -> "87-40="
$t.Cell(6, 3).Range.Text = "87-40="
# row 6, col 4: "12-12=" -> "75-27="
$t.Cell(6, 4).Range.Text = "75-27="
# row 6, col 5: "15+70=" -> "10+53="
$t.Cell(6, 5).Range.Text = "10+53="
# row 7, col 1: "3+28=" -> "79-74="
$t.Cell(7, 1).Range.Text = "79-74="
# row 7, col 2: "18+19=" -> "22+57="
$t.Cell(7, 2).Range.Text = "22+57="
# row 7, col 3: "34+54=" -> "29+58="
$t.Cell(7, 3).Range.Text = "29+58="
# row 7, col 4: "91-66=" -> "67-15="
$t.Cell(7, 4).Range.Text = "67-15="
# row 7, col 5: "39+36=" -> "75-9="
$t.Cell(7, 5).Range.Text = "75-9="
# row 8, col 1: "14+76=" -> "28-6="
$t.Cell(8, 1).Range.Text = "28-6="
# row 8, col 2: "21+48=" -> "74-27="
$t.Cell(8, 2).Range.Text = "74-27="
# row 8, col 3: "14+80=" -> "67-64="
$t.Cell(8, 3).Range.Text = "67-64="
# row 8, col 4: "4-4=" -> "70-16="
$t.Cell(8, 4).Range.Text = "70-16="
# row 8, col 5: "76-26=" -> "67-16="
$t.Cell(8, 5).Range.Text = "67-16="
# row 9, col 1: "29+8=" -> "23+19="
$t.Cell(9, 1).Range.Text = "23+19="
# row 9, col 2: "20-8=" -> "12+72="
$t.Cell(9, 2).Range.Text = "12+72="
# row 9, col 3: "30-8=" -> "81-43="
$t.Cell(9, 3).Range.Text = "81-43="
# row 9, col 4: "38+17=" -> "85-7="
$t.Cell(9, 4).Range.Text = "85-7="
# row 9, col 5: "51-42=" -> "49+50="
$t.Cell(9, 5).Range.Text = "49+50="
# row 10, col 1: "1+46=" -> "64-29="
$t.Cell(10, 1).Range.Text = "64-29="
# row 10, col 2: "81-8=" -> "18-12="
$t.Cell(10, 2).Range.Text = "18-12="
# row 10, col 3: "62+5=" -> "11+42="
$t.Cell(10, 3).Range.Text = "11+42="
# row 10, col 4: "34+48=" -> "81-18="
$t.Cell(10, 4).Range.Text = "81-18="
# row 10, col 5: "34-0=" -> "64+1="
$t.Cell(10, 5).Range.Text = "64+1="
# row 11, col 1: "86-9=" -> "76+20="
$t.Cell(11, 1).Range.Text = "76+20="
# row 11, col 2: "66-47=" -> "72-3="
$t.Cell(11, 2).Range.Text = "72-3="
# row 11, col 3: "88-79=" -> "68-32="
$t.Cell(11, 3).Range.Text = "68-32="
# row 11, col 4: "83-72=" -> "52+37="
$t.Cell(11, 4).Range.Text = "52+37="
# row 11, col 5: "90-74=" -> "91-11="
$t.Cell(11, 5).Range.Text = "91-11="
# row 12, col 1: "88-86=" -> "90-13="
$t.Cell(12, 1).Range.Text = "90-13="
# row 12, col 2: "74-6=" -> "65-16="
$t.Cell(12, 2).Range.Text = "65-16="
# row 12, col 3: "14+59=" -> "31-11="
$t.Cell(12, 3).Range.Text = "31-11="
# row 12, col 4: "0+71=" -> "97-42="
$t.Cell(12, 4).Range.Text = "97-42="
# row 12, col 5: "0+72=" -> "9+86="
$t.Cell(12, 5).Range.Text = "9+86="
# row 13, col 1: "24+31=" -> "6+5="
$t.Cell(13, 1).Range.Text = "6+5="
# row 13, col 2: "73+15=" -> "58-0="
$t.Cell(13, 2).Range.Text = "58-0="
# row 13, col 3: "41-24=" -> "21+63="
$t.Cell(13, 3).Range.Text = "21+63="
# row 13, col 4: "10+11=" -> "34+47="
$t.Cell(13, 4).Range.Text = "34+47="
# row 13, col 5: "31+32=" -> "44-27="
$t.Cell(13, 5).Range.Text = "44-27="
# row 14, col 1: "72+7=" -> "44-39="
$t.Cell(14, 1).Range.Text = "44-39="
# row 14, col 2: "83+14=" -> "78+17="
$t.Cell(14, 2).Range.Text = "78+17="
# row 14, col 3: "63+30=" -> "78-47="
$t.Cell(14, 3).Range.Text = "78-47="
# row 14, col 4: "82-8=" -> "48+26="
$t.Cell(14, 4).Range.Text = "48+26="
# row 14, col 5: "76-29=" -> "98-2="
$t.Cell(14, 5).Range.Text = "98-2="
# row 15, col 1: "77-22=" -> "75+5="
$t.Cell(15, 1).Range.Text = "75+5="
# row 15, col 2: "57-30=" -> "71+6="
$t.Cell(15, 2).Range.Text = "71+6="
# row 15, col 3: "27-2=" -> "55-32="
$t.Cell(15, 3).Range.Text = "55-32="
# row 15, col 4: "82-39=" -> "40-37="
$t.Cell(15, 4).Range.Text = "40-37="
# row 15, col 5: "97-92=" -> "7+74="
$t.Cell(15, 5).Range.Text = "7+74="
# row 16, col 1: "93-73=" -> "35-0="
$t.Cell(16, 1).Range.Text = "35-0="
# row 16, col 2: "25-1=" -> "44+34="
$t.Cell(16, 2).Range.Text = "44+34="
# row 16, col 3: "62-22=" -> "60-58="
$t.Cell(16, 3).Range.Text = "60-58="
# row 16, col 4: "21-21=" -> "72+14="
$t.Cell(16, 4).Range.Text = "72+14="
# row 16, col 5: "91-65=" -> "32-2="
$t.Cell(16, 5).Range.Text = "32-2="
# row 17, col 1: "0+36=" -> "80-51="
$t.Cell(17, 1).Range.Text = "80-51="
# row 17, col 2: "10+43=" -> "46-26="
$t.Cell(17, 2).Range.Text = "46-26="
# row 17, col 3: "23+11=" -> "12+47="
$t.Cell(17, 3).Range.Text = "12+47="
# row 17, col 4: "20+79=" -> "42+57="
$t.Cell(17, 4).Range.Text = "42+57="
# row 17, col 5: "85-72=" -> "45+15="
$t.Cell(17, 5).Range.Text = "45+15="
# row 18, col 1: "23+51=" -> "12+7="
$t.Cell(18, 1).Range.Text = "12+7="
# row 18, col 2: "87-10=" -> "45+15="
$t.Cell(18, 2).Range.Text = "45+15="
# row 18, col 3: "29-2=" -> "32+60="
$t.Cell(18, 3).Range.Text = "32+60="
# row 18, col 4: "66-11=" -> "77+7="
$t.Cell(18, 4).Range.Text = "77+7="
# row 18, col 5: "17+17=" -> "15+44="
$t.Cell(18, 5).Range.Text = "15+44="
# row 19, col 1: "64+16=" -> "60-32="
$t.Cell(19, 1).Range.Text = "60-32="
# row 19, col 2: "17-14=" -> "29+63="
$t.Cell(19, 2).Range.Text = "29+63="
# row 19, col 3: "6+70=" -> "54-13="
$t.Cell(19, 3).Range.Text = "54-13="
# row 19, col 4: "73-41=" -> "13+18="
$t.Cell(19, 4).Range.Text = "13+18="
# row 19, col 5: "46+28=" -> "44+55="
$t.Cell(19, 5).Range.Text = "44+55="
# row 20, col 1: "81-8=" -> "96-34="
$t.Cell(20, 1).Range.Text = "96-34="
# row 20, col 2: "88-64=" -> "85-35="
$t.Cell(20, 2).Range.Text = "85-35="
# row 20, col 3: "10+21=" -> "40+39="
$t.Cell(20, 3).Range.Text = "40+39="
# row 20, col 4: "2+67=" -> "51+19="
$t.Cell(20, 4).Range.Text = "51+19="
# row 20, col 5: "91-11=" -> "98-84="
$t.Cell(20, 5).Range.Text = "98-84="
